$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 for a new "Black Rock" exposure site.
# This pushes the existing Brighton/Doveton/Glen Waverley/Mordialloc rows down by one.
$ws.Range("A4").EntireRow.Insert()
Write-Host "Inserted row 4"

$ws.Cells.Item(4, 1).Value = "Black Rock"
$ws.Cells.Item(4, 2).Value = "Smile Buffalo Thai restaurant  305 Beach Road, Black Rock VIC 3193"
$ws.Cells.Item(4, 3).Value = "27/12/20 7.30pm - 9pm"
$ws.Cells.Item(4, 4).Value = "Case dined for dinner"

# Insert two new rows before the second Mordialloc row (now at row 8) for two new
# "Melbourne" exposure sites.
$ws.Range("A8:A9").EntireRow.Insert()
Write-Host "Inserted rows 8:9"

$ws.Cells.Item(8, 1).Value = "Melbourne"
$ws.Cells.Item(8, 2).Value = "Left Bank Melbourne, 1 Southbank Blvd"
$ws.Cells.Item(8, 3).Value = "25/12/20 12pm - 3pm"
$ws.Cells.Item(8, 4).Value = "Case ate in store"

$ws.Cells.Item(9, 1).Value = "Melbourne"
$ws.Cells.Item(9, 2).Value = "Melbourne Central Lion Hotel, 211 La Trobe Street"
$ws.Cells.Item(9, 3).Value = "28/12/2020 10pm - 12.30am"
$ws.Cells.Item(9, 4).Value = "Case attended Venue"

Write-Host "Done filling new rows"
